# Apply cryptos list update (Thu May  2 18:11:19 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.499.70"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "3.003.69"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'564.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'139.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.75%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "'0.523"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("D9").Value = "2.987.58"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("D11").Value = "'5.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.15%  "
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").Value = "'0.0000231"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").Value = "'33.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "3.497.85"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "'7.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.01%  "
$ws.Range("D18").Value = "2.999.15"
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("D19").Value = "59.398.05"
$ws.Range("E19").Value = "  +2.56%  "
$ws.Range("D20").Value = "'431.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").Value = "'13.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("D22").Value = "'0.718"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("D23").Value = "'7.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").Value = "'13.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("D25").Value = "'80.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'2.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.84%  "
$ws.Range("D29").Value = "'2.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("D30").Value = "'7.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").Value = "'25.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").Value = "'6.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "'0.0993"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("D34").Value = "0.0₃0781"
$ws.Range("E34").Value = "  +17.37%  "
$ws.Range("D35").Value = "'0.993"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.95%  "
$ws.Range("E36").Value = "  +4.38%  "
$ws.Range("D37").Value = "'2.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'49.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").Value = "'8.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.86%  "
$ws.Range("D40").Value = "'2.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.60%  "
$ws.Range("D41").Value = "'404.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.81%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0354"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.781.82"
$ws.Range("E43").Value = "  +5.25%  "
$ws.Range("D44").Value = "'0.109"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  +6.92%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'34.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.75%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'121.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.111"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("B50").Value = "Fetch.AI"
$ws.Range("C50").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D50").Value = "'2.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").Value = "'23.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.75%  "

Write-Host "Applied all updates"
